$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.693.56'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '1.685.19'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.40'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '30.43'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.72%  '
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('E10').Value = '  +2.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0906'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = '1.927.17'
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('E13').Value = '  +10.91%  '
$ws.Range('D14').Value = '1.695.25'
$ws.Range('E14').Value = '  +3.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.619'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.98'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('D17').Value = '30.691.00'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.31'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.08'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').Value = '0.0₃0714'
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.20'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.28'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.46'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.84'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0498'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.48'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.50%  '
$ws.Range('D33').Value = '1.509.15'
$ws.Range('E33').Value = '  +4.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.28'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.74'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.46%  '
$ws.Range('B36').Value = 'Aave'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '84.16'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +8.89%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.03'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('E38').Value = '  +3.75%  '
$ws.Range('E39').Value = '  -5.01%  '
$ws.Range('E40').Value = '  +4.33%  '
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.839'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.99'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.23%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0500'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.69%  '
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '51.61'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -6.08%  '
$ws.Range('D48').Value = '1.819.59'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.42'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '94.76'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.00%  '
$ws.Range('D51').Value = '0.0₆0113'
$ws.Range('E51').Value = '  -1.20%  '
